$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the (empty) first row, shifting everything up by one.
$ws.Rows.Item(1).Delete()

# Re-fit column M's width to its (now-changed) content.
$ws.Columns.Item(13).AutoFit()

# Select a cell similar to what the author ended up with.
$ws.Range("E31").Select()
